$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column (D) values remain stored as text, matching the
# source data (prices are formatted strings, not numeric values).
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '29.376.14'
$ws.Range("E2").Value = '  +0.09%  '

$ws.Range("D3").Value = '1.841.56'
$ws.Range("E3").Value = '  -0.20%  '

$ws.Range("E4").Value = '  +0.12%  '

$ws.Range("D5").Value = '238.96'
$ws.Range("E5").Value = '  -0.41%  '

$ws.Range("D6").Value = '0.6266'
$ws.Range("E6").Value = '  -0.13%  '

$ws.Range("E7").Value = '  +0.13%  '

$ws.Range("D8").Value = '0.07431'
$ws.Range("E8").Value = '  -0.78%  '

$ws.Range("E9").Value = '  -0.15%  '

$ws.Range("E10").Value = '  +1.86%  '

$ws.Range("D11").Value = '0.07725'
$ws.Range("E11").Value = '  -0.15%  '

$ws.Range("D12").Value = '1.845.46'
$ws.Range("E12").Value = '  +0.01%  '

$ws.Range("D13").Value = '4.973'
$ws.Range("E13").Value = '  -0.24%  '

$ws.Range("D14").Value = '0.6742'
$ws.Range("E14").Value = '  -0.83%  '

$ws.Range("E15").Value = '  -2.35%  '

$ws.Range("D16").Value = '81.73'
$ws.Range("E16").Value = '  -0.32%  '

$ws.Range("D17").Value = '6.207'
$ws.Range("E17").Value = '  +0.61%  '

$ws.Range("D18").Value = '29.425.35'
$ws.Range("E18").Value = '  +0.15%  '

$ws.Range("D19").Value = '233.51'
$ws.Range("E19").Value = '  +1.95%  '

$ws.Range("E20").Value = '  -0.10%  '

$ws.Range("E21").Value = '  +0.17%  '

$ws.Range("D22").Value = '7.294'
$ws.Range("E22").Value = '  -2.80%  '

$ws.Range("D23").Value = '1.000'
$ws.Range("E23").Value = '  +0.12%  '

$ws.Range("D24").Value = '157.93'
$ws.Range("E24").Value = '  -0.36%  '

$ws.Range("D25").Value = '8.502'
$ws.Range("E25").Value = '  +0.87%  '

$ws.Range("E26").Value = '  -1.80%  '

$ws.Range("D27").Value = '17.31'
$ws.Range("E27").Value = '  -1.25%  '

$ws.Range("D28").Value = '0.07222'
$ws.Range("E28").Value = '  +10.64%  '

$ws.Range("D29").Value = '1.466'
$ws.Range("E29").Value = '  +3.91%  '

$ws.Range("D30").Value = '1.480'
$ws.Range("E30").Value = '  +0.15%  '

$ws.Range("D31").Value = '4.040'
$ws.Range("E31").Value = '  -1.63%  '

$ws.Range("D32").Value = '4.031'
$ws.Range("E32").Value = '  -1.48%  '

$ws.Range("E33").Value = '  -0.69%  '

$ws.Range("D34").Value = '1.140'
$ws.Range("E34").Value = '  -0.12%  '

$ws.Range("D35").Value = '0.6967'
$ws.Range("E35").Value = '  +0.06%  '

$ws.Range("E36").Value = '  -0.15%  '

$ws.Range("E37").Value = '  +0.44%  '

$ws.Range("D38").Value = '6.931'
$ws.Range("E38").Value = '  +2.43%  '

$ws.Range("D39").Value = '2.815'
$ws.Range("E39").Value = '  -0.74%  '

$ws.Range("D40").Value = '1.233.62'
$ws.Range("E40").Value = '  -2.26%  '

$ws.Range("D41").Value = '0.9621'
$ws.Range("E41").Value = '  +4.59%  '

$ws.Range("D42").Value = '1.001'
$ws.Range("E42").Value = '  +0.19%  '

$ws.Range("D43").Value = '1.999.30'
$ws.Range("E43").Value = '  -0.46%  '

$ws.Range("D44").Value = '100.87'
$ws.Range("E44").Value = '  -0.38%  '

$ws.Range("D45").Value = '65.38'
$ws.Range("E45").Value = '  -1.20%  '

$ws.Range("B46").Value = 'RenderToken'
$ws.Range("C46").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D46").Value = '1.716'
$ws.Range("E46").Value = '  -0.57%  '

$ws.Range("B47").Value = 'Aptos'
$ws.Range("C47").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D47").Value = '6.936'
$ws.Range("E47").Value = '  -2.01%  '

$ws.Range("B48").Value = 'TheSandbox'
$ws.Range("C48").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D48").Value = '0.3899'
$ws.Range("E48").Value = '  -1.43%  '

$ws.Range("B49").Value = 'Algorand'
$ws.Range("C49").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D49").Value = '0.1131'
$ws.Range("E49").Value = '  -2.63%  '

$ws.Range("B50").Value = 'EnergySwap'
$ws.Range("C50").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D50").Value = '8.838'
$ws.Range("E50").Value = '  -1.46%  '

$ws.Range("B51").Value = 'Cronos'
$ws.Range("C51").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D51").Value = '0.05658'
$ws.Range("E51").Value = '  -0.65%  '

